$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 corresponds to the "Undecided" category.
# Update the daily counts (B9:F9) and the weekly total (G9).
$ws.Range("B9").Value = 12
$ws.Range("C9").Value = 9
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 12
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = 55
